$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - header labels (date range text)
$ws.Range("D3").Value = "OCT' 21"
$ws.Range("E3").Value = "OCT' 22"
$ws.Range("F3").Value = "OCT' 22"

# Row 4 - PU1
$ws.Range("D4").Value = 1507.04
$ws.Range("E4").Value = 1535.17
$ws.Range("F4").Value = 1528.7

# Row 5 - PU10
$ws.Range("D5").Value = 131.63
$ws.Range("E5").Value = 139.05
$ws.Range("F5").Value = 169.97

# Row 6 - PU11
$ws.Range("D6").Value = 0.99
$ws.Range("E6").Value = 2.56
$ws.Range("F6").Value = 9.88

# Row 7 - PU15
$ws.Range("D7").Value = 4.04
$ws.Range("E7").Value = 5.23
$ws.Range("F7").Value = 5.5

# Row 8 - PU32
$ws.Range("C8").Value = 305.42
$ws.Range("D8").Value = 209.56
$ws.Range("E8").Value = 173.34
$ws.Range("F8").Value = 230.28
